$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.080.61'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.31'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5226'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3758'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07265'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8995'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08177'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.49%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.928.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '96.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.293'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008573'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.110.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.078'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.69'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.407'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.289'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.732'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.96'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.785'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('E30').Value = '  -2.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09224'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05036'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7886'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.210'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.434'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.74%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.975'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.600'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5723'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01984'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.074'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.040'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.554'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '116.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1516'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4864'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.05'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.624'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '38.15'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '63.52'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05927'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.04%  '
